$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.854.74'
$ws.Range('E2').Value = '  -0.79%  '
$ws.Range('D3').Value = '3.515.19'
$ws.Range('E3').Value = '  +2.31%  '
$ws.Range('E4').Value = '  +0.08%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '592.64'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +0.91%  '
$ws.Range('E6').Value = '  -1.32%  '
$ws.Range('D7').Value = '3.519.64'
$ws.Range('E7').Value = '  +2.47%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('E10').Value = '  +0.47%  '
$ws.Range('E11').Value = '  -3.96%  '
$ws.Range('E12').Value = '  +1.97%  '
$ws.Range('D13').Value = '4.112.87'
$ws.Range('E13').Value = '  +2.36%  '
$ws.Range('D14').Value = '3.527.92'
$ws.Range('E14').Value = '  +0.63%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '27.00'
$c.Style = "Normal"
$ws.Range('E15').Value = '  +1.83%  '
$ws.Range('E16').Value = '  +0.51%  '
$ws.Range('E17').Value = '  +1.34%  '
$ws.Range('D18').Value = '64.899.67'
$ws.Range('E18').Value = '  -0.47%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '10.01'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +3.38%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '5.80'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -0.34%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '14.12'
$c.Style = "Normal"
$ws.Range('E21').Value = '  +4.42%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '387.30'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +0.45%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '0.569'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +2.88%  '
$ws.Range('D24').Value = '3.665.81'
$ws.Range('E24').Value = '  +2.67%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '73.78'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +1.60%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('E27').Value = '  +3.88%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '7.59'
$c.Style = "Normal"
$ws.Range('E28').Value = '  +5.83%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '0.996'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -0.17%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '2.24'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +1.47%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '8.15'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +0.21%  '
$ws.Range('D32').Value = '3.532.01'
$ws.Range('E32').Value = '  +2.57%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '23.68'
$c.Style = "Normal"
$ws.Range('E34').Value = '  +2.27%  '
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '1.35'
$c.Style = "Normal"
$ws.Range('E35').Value = '  +14.10%  '
$ws.Range('E36').Value = '  -0.72%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '168.94'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +0.13%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '1.55'
$c.Style = "Normal"
$ws.Range('E38').Value = '  +6.22%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '6.82'
$c.Style = "Normal"
$ws.Range('E39').Value = '  +0.68%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '4.91'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +4.91%  '
$ws.Range('E41').Value = '  +5.01%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.817'
$c.Style = "Normal"
$ws.Range('E42').Value = '  +0.48%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '26.13'
$c.Style = "Normal"
$ws.Range('E43').Value = '  +14.32%  '
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('E45').Value = '  -0.47%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '4.39'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +0.85%  '
$ws.Range('E47').Value = '  +5.10%  '
$ws.Range('E48').Value = '  +1.76%  '
$ws.Range('E49').Value = '  +5.00%  '
$ws.Range('D50').Value = '2.382.10'
$ws.Range('E50').Value = '  +9.95%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '302.46'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +11.76%  '
